# Rename the "params" sheet to "decomps"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "decomps"

# Update the view state: scroll/show row 2 at the top and select L18
# (mirrors topLeftCell="A2" + <selection activeCell="L18" sqref="L18"/>)
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("L18").Select()
